$wb = $excel.ActiveWorkbook

$wsSurvey   = $wb.Worksheets.Item("survey")
$wsChoices  = $wb.Worksheets.Item("choices")
$wsSettings = $wb.Worksheets.Item("settings")

# Update the value list referenced by the plant_type question from "plants" to "planting".
$wsSurvey.Range("D5").Value = "planting"

# Restore the selections on each sheet (settings and choices keep their own selection,
# survey becomes the active/selected tab last so it ends up the active sheet).
[void]$wsSettings.Range("B15").Select()
[void]$wsChoices.Range("C6").Select()
[void]$wsSurvey.Range("D5").Select()
